$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("opcionPrincipal" / "Clientes") ---
$ws.Range("C1").Value = "opcionPrincipal"
$ws.Range("C2").Value = "Clientes"

# --- Column D ("Tip_documento" / "Cedula de ciudadania") ---
$ws.Range("D1").Value = "Tip_documento"
$ws.Range("D2").Value = "Cédula de ciudadanía"

# Copy down the existing formatting (border/fill/alignment) from column B
# onto the two new columns so the new cells match the sheet's existing
# look - header row, data row and the empty trailing row.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null

$ws.Range("B3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null

# --- Column widths ---
$ws.Range("B1").ColumnWidth = 25
$ws.Range("C1").ColumnWidth = 24.166666666666668
$ws.Range("D1").ColumnWidth = 30

# --- Selection moves to the newly added D2 cell ---
$ws.Range("D2").Select() | Out-Null
